$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 15.57800190790504
$ws.Cells.Item(2, 3).Value = 8.00005591065424
$ws.Cells.Item(2, 4).Value = 7.940792452307045
$ws.Cells.Item(2, 5).Value = 12.99624494455836
$ws.Cells.Item(2, 6).Value = 38.13641248422856
$ws.Cells.Item(2, 8).Value = 7.344005520526261
$ws.Cells.Item(2, 10).Value = 10.24051627375567
$ws.Cells.Item(2, 11).Value = 11.68610621854213
$ws.Cells.Item(2, 12).Value = 10.79260582059076
$ws.Cells.Item(2, 14).Value = 21.20514023523401
$ws.Cells.Item(2, 15).Value = 29.45693456936715
$ws.Cells.Item(3, 2).Value = 15.40309458605463
$ws.Cells.Item(3, 3).Value = 7.979124672058223
$ws.Cells.Item(3, 4).Value = 7.925965987496923
$ws.Cells.Item(3, 5).Value = 13.01381980262528
$ws.Cells.Item(3, 6).Value = 38.21399427253784
$ws.Cells.Item(3, 8).Value = 7.344005520526261
$ws.Cells.Item(3, 10).Value = 10.260062665813
$ws.Cells.Item(3, 11).Value = 11.55980489088758
$ws.Cells.Item(3, 12).Value = 10.79251960207046
$ws.Cells.Item(3, 14).Value = 21.2658948785684
$ws.Cells.Item(3, 15).Value = 29.53561375106047
$ws.Cells.Item(4, 2).Value = 15.29757681231026
$ws.Cells.Item(4, 3).Value = 7.966175663296478
$ws.Cells.Item(4, 4).Value = 7.917934746495056
$ws.Cells.Item(4, 5).Value = 13.02618630448234
$ws.Cells.Item(4, 6).Value = 38.26871890067806
$ws.Cells.Item(4, 8).Value = 7.344005520526261
$ws.Cells.Item(4, 10).Value = 10.2728155832853
$ws.Cells.Item(4, 11).Value = 11.48348350119372
$ws.Cells.Item(4, 12).Value = 10.79376169256765
$ws.Cells.Item(4, 14).Value = 21.30495113743881
$ws.Cells.Item(4, 15).Value = 29.58887585530297
$ws.Cells.Item(5, 2).Value = 15.25509596058929
$ws.Cells.Item(5, 3).Value = 7.960875547451754
$ws.Cells.Item(5, 4).Value = 7.914934082723714
$ws.Cells.Item(5, 5).Value = 13.03162245094344
$ws.Cells.Item(5, 6).Value = 38.29280074596394
$ws.Cells.Item(5, 8).Value = 7.344005520526261
$ws.Cells.Item(5, 10).Value = 10.27820188468305
$ws.Cells.Item(5, 11).Value = 11.45272296537465
$ws.Cells.Item(5, 12).Value = 10.79459463752374
$ws.Cells.Item(5, 14).Value = 21.32130888670463
$ws.Cells.Item(5, 15).Value = 29.6118252015654
$ws.Cells.Item(6, 2).Value = 15.24807471192453
$ws.Cells.Item(6, 3).Value = 7.959994093636411
$ws.Cells.Item(6, 4).Value = 7.914452330602828
$ws.Cells.Item(6, 5).Value = 13.03254909328339
$ws.Cells.Item(6, 6).Value = 38.2969070439066
$ws.Cells.Item(6, 8).Value = 7.344005520526261
$ws.Cells.Item(6, 10).Value = 10.27910772751406
$ws.Cells.Item(6, 11).Value = 11.4476367305256
$ws.Cells.Item(6, 12).Value = 10.79475271581458
$ws.Cells.Item(6, 14).Value = 21.32405181245001
$ws.Cells.Item(6, 15).Value = 29.61571106933754
$ws.Cells.Item(7, 2).Value = 15.29700174021526
$ws.Cells.Item(7, 3).Value = 7.966104276883248
$ws.Cells.Item(7, 4).Value = 7.917893173458085
$ws.Cells.Item(7, 5).Value = 13.02625801136185
$ws.Cells.Item(7, 6).Value = 38.26903646738687
$ws.Cells.Item(7, 8).Value = 7.344005520526261
$ws.Cells.Item(7, 10).Value = 10.27288745742099
$ws.Cells.Item(7, 11).Value = 11.48306723081732
$ws.Cells.Item(7, 12).Value = 10.7937716013813
$ws.Cells.Item(7, 14).Value = 21.30516995239519
$ws.Cells.Item(7, 15).Value = 29.58918031997287
$ws.Cells.Item(8, 2).Value = 15.51733279963872
$ws.Cells.Item(8, 3).Value = 7.992859664414878
$ws.Cells.Item(8, 4).Value = 7.935459288181872
$ws.Cells.Item(8, 5).Value = 13.00197809643556
$ws.Cells.Item(8, 6).Value = 38.16169052082187
$ws.Cells.Item(8, 8).Value = 7.344005520526261
$ws.Cells.Item(8, 10).Value = 10.24710019703286
$ws.Cells.Item(8, 11).Value = 11.6423221195092
$ws.Cells.Item(8, 12).Value = 10.7923080151782
$ws.Cells.Item(8, 14).Value = 21.22572544864267
$ws.Cells.Item(8, 15).Value = 29.48303473815549
$ws.Cells.Item(9, 2).Value = 15.96219757783828
$ws.Cells.Item(9, 3).Value = 8.044507118706548
$ws.Cells.Item(9, 4).Value = 7.978297886676544
$ws.Cells.Item(9, 5).Value = 12.96684177193301
$ws.Cells.Item(9, 6).Value = 38.00749257953449
$ws.Cells.Item(9, 8).Value = 7.344005520526261
$ws.Cells.Item(9, 10).Value = 10.20247279943892
$ws.Cells.Item(9, 11).Value = 11.96292707200528
$ws.Cells.Item(9, 12).Value = 10.79966095780771
$ws.Cells.Item(9, 14).Value = 21.08378191923779
$ws.Cells.Item(9, 15).Value = 29.31421605821443
$ws.Cells.Item(10, 2).Value = 16.29401551694193
$ws.Cells.Item(10, 3).Value = 8.081888523913694
$ws.Cells.Item(10, 4).Value = 8.01472155480095
$ws.Cells.Item(10, 5).Value = 12.94859958686866
$ws.Cells.Item(10, 6).Value = 37.92859643713278
$ws.Cells.Item(10, 8).Value = 7.344005520526261
$ws.Cells.Item(10, 10).Value = 10.17327893401039
$ws.Cells.Item(10, 11).Value = 12.20159610485681
$ws.Cells.Item(10, 12).Value = 10.81122133576604
$ws.Cells.Item(10, 14).Value = 20.98785243877107
$ws.Cells.Item(10, 15).Value = 29.21420772065378
$ws.Cells.Item(11, 2).Value = 16.44544881833606
$ws.Cells.Item(11, 3).Value = 8.098758029753947
$ws.Cells.Item(11, 4).Value = 8.032327841286458
$ws.Cells.Item(11, 5).Value = 12.94193712707575
$ws.Cells.Item(11, 6).Value = 37.90018083472007
$ws.Cells.Item(11, 8).Value = 7.344005520526261
$ws.Cells.Item(11, 10).Value = 10.16077232783972
$ws.Cells.Item(11, 11).Value = 12.31043777752737
$ws.Cells.Item(11, 12).Value = 10.81780003650477
$ws.Cells.Item(11, 14).Value = 20.94600800697846
$ws.Cells.Item(11, 15).Value = 29.17393582521677
$ws.Cells.Item(12, 2).Value = 16.5028132167795
$ws.Cells.Item(12, 3).Value = 8.105125354075746
$ws.Cells.Item(12, 4).Value = 8.039140452357763
$ws.Cells.Item(12, 5).Value = 12.93964873636826
$ws.Cells.Item(12, 6).Value = 37.89049548342528
$ws.Cells.Item(12, 8).Value = 7.344005520526261
$ws.Cells.Item(12, 10).Value = 10.15614722683079
$ws.Cells.Item(12, 11).Value = 12.35165793605704
$ws.Cells.Item(12, 12).Value = 10.82047933528633
$ws.Cells.Item(12, 14).Value = 20.93041929946628
$ws.Cells.Item(12, 15).Value = 29.15943743722778
$ws.Cells.Item(13, 2).Value = 16.49045867286317
$ws.Cells.Item(13, 3).Value = 8.103754982502748
$ws.Cells.Item(13, 4).Value = 8.037666823172014
$ws.Cells.Item(13, 5).Value = 12.94013116301634
$ws.Cells.Item(13, 6).Value = 37.89253358631603
$ws.Cells.Item(13, 8).Value = 7.344005520526261
$ws.Cells.Item(13, 10).Value = 10.15713839996863
$ws.Cells.Item(13, 11).Value = 12.34278079694473
$ws.Cells.Item(13, 12).Value = 10.81989396323709
$ws.Cells.Item(13, 14).Value = 20.93376520047165
$ws.Cells.Item(13, 15).Value = 29.16252648160068
$ws.Cells.Item(14, 2).Value = 16.45016810784008
$ws.Cells.Item(14, 3).Value = 8.09928229757981
$ws.Cells.Item(14, 4).Value = 8.032885424310377
$ws.Cells.Item(14, 5).Value = 12.94174416357167
$ws.Cells.Item(14, 6).Value = 37.89936247029572
$ws.Cells.Item(14, 8).Value = 7.344005520526261
$ws.Cells.Item(14, 10).Value = 10.16038959806534
$ws.Cells.Item(14, 11).Value = 12.31382909740159
$ws.Cells.Item(14, 12).Value = 10.81801670456395
$ws.Cells.Item(14, 14).Value = 20.94472037315113
$ws.Cells.Item(14, 15).Value = 29.17272796498211
$ws.Cells.Item(15, 2).Value = 16.42549008805199
$ws.Cells.Item(15, 3).Value = 8.096539902251802
$ws.Cells.Item(15, 4).Value = 8.029975516712899
$ws.Cells.Item(15, 5).Value = 12.94276269471382
$ws.Cells.Item(15, 6).Value = 37.90368535742229
$ws.Cells.Item(15, 8).Value = 7.344005520526261
$ws.Cells.Item(15, 10).Value = 10.16239547884151
$ws.Cells.Item(15, 11).Value = 12.29609487130356
$ws.Cells.Item(15, 12).Value = 10.81689127181918
$ws.Cells.Item(15, 14).Value = 20.95146415071346
$ws.Cells.Item(15, 15).Value = 29.1790745823246
$ws.Cells.Item(16, 2).Value = 16.2841244947628
$ws.Cells.Item(16, 3).Value = 8.080783224074947
$ws.Cells.Item(16, 4).Value = 8.013591502972874
$ws.Cells.Item(16, 5).Value = 12.94906785461174
$ws.Cells.Item(16, 6).Value = 37.93060389575074
$ws.Cells.Item(16, 8).Value = 7.344005520526261
$ws.Cells.Item(16, 10).Value = 10.17411180758335
$ws.Cells.Item(16, 11).Value = 12.19448550423448
$ws.Cells.Item(16, 12).Value = 10.81081781721943
$ws.Cells.Item(16, 14).Value = 20.99062307891045
$ws.Cells.Item(16, 15).Value = 29.21694473631469
$ws.Cells.Item(17, 2).Value = 16.197490191035
$ws.Cells.Item(17, 3).Value = 8.071081383499752
$ws.Cells.Item(17, 4).Value = 8.003803507023676
$ws.Cells.Item(17, 5).Value = 12.95335439002678
$ws.Cells.Item(17, 6).Value = 37.94903222165441
$ws.Cells.Item(17, 8).Value = 7.344005520526261
$ws.Cells.Item(17, 10).Value = 10.18149730871435
$ws.Cells.Item(17, 11).Value = 12.13219577694938
$ws.Cells.Item(17, 12).Value = 10.80742880871963
$ws.Cells.Item(17, 14).Value = 21.0151045671552
$ws.Cells.Item(17, 15).Value = 29.24151498007375
$ws.Cells.Item(18, 2).Value = 16.14770929310227
$ws.Cells.Item(18, 3).Value = 8.065488468056852
$ws.Cells.Item(18, 4).Value = 7.998271542692365
$ws.Cells.Item(18, 5).Value = 12.95597392587012
$ws.Cells.Item(18, 6).Value = 37.96033526843907
$ws.Cells.Item(18, 8).Value = 7.344005520526261
$ws.Cells.Item(18, 10).Value = 10.18581811357348
$ws.Cells.Item(18, 11).Value = 12.09639582462701
$ws.Cells.Item(18, 12).Value = 10.80560387699677
$ws.Cells.Item(18, 14).Value = 21.02935463477075
$ws.Cells.Item(18, 15).Value = 29.25613865344943
$ws.Cells.Item(19, 2).Value = 16.1308642723279
$ws.Cells.Item(19, 3).Value = 8.063592666060069
$ws.Cells.Item(19, 4).Value = 7.996415425240358
$ws.Cells.Item(19, 5).Value = 12.95688733226648
$ws.Cells.Item(19, 6).Value = 37.96428311201253
$ws.Cells.Item(19, 8).Value = 7.344005520526261
$ws.Cells.Item(19, 10).Value = 10.18729358977139
$ws.Cells.Item(19, 11).Value = 12.08428035467542
$ws.Cells.Item(19, 12).Value = 10.80500739162369
$ws.Cells.Item(19, 14).Value = 21.03420851784796
$ws.Cells.Item(19, 15).Value = 29.26117437308656
$ws.Cells.Item(20, 2).Value = 16.20670787001573
$ws.Cells.Item(20, 3).Value = 8.072115482092062
$ws.Cells.Item(20, 4).Value = 8.004835356417486
$ws.Cells.Item(20, 5).Value = 12.95288214411783
$ws.Cells.Item(20, 6).Value = 37.94699767655101
$ws.Cells.Item(20, 8).Value = 7.344005520526261
$ws.Cells.Item(20, 10).Value = 10.18070357166511
$ws.Cells.Item(20, 11).Value = 12.13882403626474
$ws.Cells.Item(20, 12).Value = 10.8077767181404
$ws.Cells.Item(20, 14).Value = 21.01248099230391
$ws.Cells.Item(20, 15).Value = 29.23884855970643
$ws.Cells.Item(21, 2).Value = 16.46200227953029
$ws.Cells.Item(21, 3).Value = 8.100596608312074
$ws.Cells.Item(21, 4).Value = 8.034285918068884
$ws.Cells.Item(21, 5).Value = 12.94126402685412
$ws.Cells.Item(21, 6).Value = 37.89732748616627
$ws.Cells.Item(21, 8).Value = 7.344005520526261
$ws.Cells.Item(21, 10).Value = 10.15943163628028
$ws.Cells.Item(21, 11).Value = 12.32233307316312
$ws.Cells.Item(21, 12).Value = 10.8185630103411
$ws.Cells.Item(21, 14).Value = 20.94149561115327
$ws.Cells.Item(21, 15).Value = 29.16971113562411
$ws.Cells.Item(22, 2).Value = 16.62894525478709
$ws.Cells.Item(22, 3).Value = 8.119088868944207
$ws.Cells.Item(22, 4).Value = 8.054379817317093
$ws.Cells.Item(22, 5).Value = 12.93503761807558
$ws.Cells.Item(22, 6).Value = 37.87113101478865
$ws.Cells.Item(22, 8).Value = 7.344005520526261
$ws.Cells.Item(22, 10).Value = 10.14617530531556
$ws.Cells.Item(22, 11).Value = 12.44227556295379
$ws.Cells.Item(22, 12).Value = 10.82670793672505
$ws.Cells.Item(22, 14).Value = 20.89659927037041
$ws.Cells.Item(22, 15).Value = 29.12890761412936
$ws.Cells.Item(23, 2).Value = 16.53985246662525
$ws.Cells.Item(23, 3).Value = 8.109230781076548
$ws.Cells.Item(23, 4).Value = 8.043579128869554
$ws.Cells.Item(23, 5).Value = 12.93823596520028
$ws.Cells.Item(23, 6).Value = 37.88453927226544
$ws.Cells.Item(23, 8).Value = 7.344005520526261
$ws.Cells.Item(23, 10).Value = 10.15319146957421
$ws.Cells.Item(23, 11).Value = 12.37827048962737
$ws.Cells.Item(23, 12).Value = 10.82226118864843
$ws.Cells.Item(23, 14).Value = 20.92042472306079
$ws.Cells.Item(23, 15).Value = 29.15028407159734
$ws.Cells.Item(24, 2).Value = 16.20254047152741
$ws.Cells.Item(24, 3).Value = 8.07164801336595
$ws.Cells.Item(24, 4).Value = 8.004368560243149
$ws.Cells.Item(24, 5).Value = 12.95309516317268
$ws.Cells.Item(24, 6).Value = 37.94791528792341
$ws.Cells.Item(24, 8).Value = 7.344005520526261
$ws.Cells.Item(24, 10).Value = 10.18106218721662
$ws.Cells.Item(24, 11).Value = 12.13582736255838
$ws.Cells.Item(24, 12).Value = 10.80761904362062
$ws.Cells.Item(24, 14).Value = 21.0136665642401
$ws.Cells.Item(24, 15).Value = 29.24005249743339
$ws.Cells.Item(25, 2).Value = 15.84077124276564
$ws.Cells.Item(25, 3).Value = 8.030629455866251
$ws.Cells.Item(25, 4).Value = 7.965826951853358
$ws.Cells.Item(25, 5).Value = 12.97501467847139
$ws.Cells.Item(25, 6).Value = 38.04317100009751
$ws.Cells.Item(25, 8).Value = 7.344005520526261
$ws.Cells.Item(25, 10).Value = 10.2139125534743
$ws.Cells.Item(25, 11).Value = 11.87550603525677
$ws.Cells.Item(25, 12).Value = 10.79658461531664
$ws.Cells.Item(25, 14).Value = 21.12070782624119
$ws.Cells.Item(25, 15).Value = 29.3556702568766
